$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates reflecting the 2020-12-28 data refresh: nombre_aides (C) and montant_total (D)
# Values are forced to remain text cells (matching the source inlineStr string typing)
# by writing as "@" (text) number format and resetting style to Normal afterwards,
# so no residual style index is left on the cell.

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1413"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "6931946.29"
$ws.Range("D3").Style = "Normal"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "625"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5752346.20"
$ws.Range("D4").Style = "Normal"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "73"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1201921.05"
$ws.Range("D6").Style = "Normal"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "90"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "340369.44"
$ws.Range("D9").Style = "Normal"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "481"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2628700.26"
$ws.Range("D10").Style = "Normal"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "64"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "852358.00"
$ws.Range("D12").Style = "Normal"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "20"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "443076.00"
$ws.Range("D13").Style = "Normal"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "31"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "99068.05"
$ws.Range("D14").Style = "Normal"

$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "175"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "831569.56"
$ws.Range("D49").Style = "Normal"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "96"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "655368.00"
$ws.Range("D50").Style = "Normal"

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "43"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "379610.00"
$ws.Range("D51").Style = "Normal"

$ws.Range("C53").NumberFormat = "@"
$ws.Range("C53").Value = "4"
$ws.Range("C53").Style = "Normal"
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = "110000.00"
$ws.Range("D53").Style = "Normal"

$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "1014"
$ws.Range("C56").Style = "Normal"
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = "5653744.71"
$ws.Range("D56").Style = "Normal"

$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = "490"
$ws.Range("C57").Style = "Normal"
$ws.Range("D57").NumberFormat = "@"
$ws.Range("D57").Value = "4395859.82"
$ws.Range("D57").Style = "Normal"

$ws.Range("C58").NumberFormat = "@"
$ws.Range("C58").Value = "181"
$ws.Range("C58").Style = "Normal"
$ws.Range("D58").NumberFormat = "@"
$ws.Range("D58").Value = "1805738.84"
$ws.Range("D58").Style = "Normal"

$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value = "14"
$ws.Range("C60").Style = "Normal"
$ws.Range("D60").NumberFormat = "@"
$ws.Range("D60").Value = "280532.77"
$ws.Range("D60").Style = "Normal"

$ws.Range("C83").NumberFormat = "@"
$ws.Range("C83").Value = "120"
$ws.Range("C83").Style = "Normal"
$ws.Range("D83").NumberFormat = "@"
$ws.Range("D83").Value = "413887.23"
$ws.Range("D83").Style = "Normal"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "507"
$ws.Range("C84").Style = "Normal"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "2510796.12"
$ws.Range("D84").Style = "Normal"

$ws.Range("C85").NumberFormat = "@"
$ws.Range("C85").Value = "212"
$ws.Range("C85").Style = "Normal"
$ws.Range("D85").NumberFormat = "@"
$ws.Range("D85").Value = "1901640.84"
$ws.Range("D85").Style = "Normal"

$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "87"
$ws.Range("C86").Style = "Normal"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "1284135.25"
$ws.Range("D86").Style = "Normal"

$ws.Range("C87").NumberFormat = "@"
$ws.Range("C87").Value = "33"
$ws.Range("C87").Style = "Normal"
$ws.Range("D87").NumberFormat = "@"
$ws.Range("D87").Value = "655737.04"
$ws.Range("D87").Style = "Normal"

$ws.Range("C91").NumberFormat = "@"
$ws.Range("C91").Value = "1096"
$ws.Range("C91").Style = "Normal"
$ws.Range("D91").NumberFormat = "@"
$ws.Range("D91").Value = "5718999.26"
$ws.Range("D91").Style = "Normal"

$ws.Range("C92").NumberFormat = "@"
$ws.Range("C92").Value = "467"
$ws.Range("C92").Style = "Normal"
$ws.Range("D92").NumberFormat = "@"
$ws.Range("D92").Value = "4225438.13"
$ws.Range("D92").Style = "Normal"

$ws.Range("C93").NumberFormat = "@"
$ws.Range("C93").Value = "189"
$ws.Range("C93").Style = "Normal"
$ws.Range("D93").NumberFormat = "@"
$ws.Range("D93").Value = "2197546.63"
$ws.Range("D93").Style = "Normal"

$ws.Range("C94").NumberFormat = "@"
$ws.Range("C94").Value = "68"
$ws.Range("C94").Style = "Normal"
$ws.Range("D94").NumberFormat = "@"
$ws.Range("D94").Value = "1291722.17"
$ws.Range("D94").Style = "Normal"

$ws.Range("C97").NumberFormat = "@"
$ws.Range("C97").Value = "390"
$ws.Range("C97").Style = "Normal"
$ws.Range("D97").NumberFormat = "@"
$ws.Range("D97").Value = "1707241.76"
$ws.Range("D97").Style = "Normal"

$ws.Range("C98").NumberFormat = "@"
$ws.Range("C98").Value = "1364"
$ws.Range("C98").Style = "Normal"
$ws.Range("D98").NumberFormat = "@"
$ws.Range("D98").Value = "6781730.90"
$ws.Range("D98").Style = "Normal"

$ws.Range("C99").NumberFormat = "@"
$ws.Range("C99").Value = "528"
$ws.Range("C99").Style = "Normal"
$ws.Range("D99").NumberFormat = "@"
$ws.Range("D99").Value = "3849894.33"
$ws.Range("D99").Style = "Normal"

$ws.Range("C100").NumberFormat = "@"
$ws.Range("C100").Value = "190"
$ws.Range("C100").Style = "Normal"
$ws.Range("D100").NumberFormat = "@"
$ws.Range("D100").Value = "2277525.50"
$ws.Range("D100").Style = "Normal"

$ws.Range("C101").NumberFormat = "@"
$ws.Range("C101").Value = "62"
$ws.Range("C101").Style = "Normal"
$ws.Range("D101").NumberFormat = "@"
$ws.Range("D101").Value = "1055125.83"
$ws.Range("D101").Style = "Normal"

$ws.Range("C102").NumberFormat = "@"
$ws.Range("C102").Value = "9"
$ws.Range("C102").Style = "Normal"
$ws.Range("D102").NumberFormat = "@"
$ws.Range("D102").Value = "290000.00"
$ws.Range("D102").Style = "Normal"

$ws.Range("C110").NumberFormat = "@"
$ws.Range("C110").Value = "443"
$ws.Range("C110").Style = "Normal"
$ws.Range("D110").NumberFormat = "@"
$ws.Range("D110").Value = "1604372.63"
$ws.Range("D110").Style = "Normal"

$ws.Range("C112").NumberFormat = "@"
$ws.Range("C112").Value = "729"
$ws.Range("C112").Style = "Normal"
$ws.Range("D112").NumberFormat = "@"
$ws.Range("D112").Value = "4891453.64"
$ws.Range("D112").Style = "Normal"

$ws.Range("C113").NumberFormat = "@"
$ws.Range("C113").Value = "247"
$ws.Range("C113").Style = "Normal"
$ws.Range("D113").NumberFormat = "@"
$ws.Range("D113").Value = "3022839.27"
$ws.Range("D113").Style = "Normal"

Write-Host "Applied 2020-12-28 data updates to" (34) "rows"
